# ------------------------------------------------------------------
# Weekly CompStat refresh: new crime data collected for report week
# 8/12/2024 - 8/18/2024 (Volume 31, Number 33).
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Range("A8").Value = "Volume 31   Number  33"
# --- Row 9 ---
$ws.Range("C9").Value = "Report Covering the Week  8/12/2024  Through  8/18/2024"
# --- Row 15 ---
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -16.666666666666
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -72.222222222222
# --- Row 16 ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 84
$ws.Range("J16").Value = 89
$ws.Range("K16").Value = -5.617977528089
$ws.Range("L16").Value = 25.373134328358
$ws.Range("M16").Value = -23.636363636363
$ws.Range("N16").Value = -83.030303030303
# --- Row 17 ---
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -87.5
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -47.826086956521
$ws.Range("I17").Value = 114
$ws.Range("J17").Value = 118
$ws.Range("K17").Value = -3.389830508474
$ws.Range("L17").Value = 11.764705882352
$ws.Range("M17").Value = 40.740740740740
$ws.Range("N17").Value = -40.314136125654
# --- Row 18 ---
$ws.Range("D14").Copy($ws.Range("C18"))
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 33
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = -45
$ws.Range("L18").Value = -58.227848101265
$ws.Range("M18").Value = -35.294117647058
$ws.Range("N18").Value = -91.624365482233
# --- Row 19 ---
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 0
$ws.Range("I19").Value = 217
$ws.Range("J19").Value = 251
$ws.Range("K19").Value = -13.545816733067
$ws.Range("L19").Value = -4.824561403508
$ws.Range("M19").Value = 26.162790697674
$ws.Range("N19").Value = -47.836538461538
# --- Row 20 ---
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 24
$ws.Range("K20").Value = -66.197183098591
$ws.Range("L20").Value = -45.454545454545
$ws.Range("M20").Value = 71.428571428571
$ws.Range("N20").Value = -90.163934426229
# --- Row 21 ---
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -36.842105263157
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = -19.230769230769
$ws.Range("I21").Value = 478
$ws.Range("J21").Value = 595
$ws.Range("K21").Value = -19.663865546218
$ws.Range("L21").Value = -9.981167608286
$ws.Range("M21").Value = 8.636363636363
$ws.Range("N21").Value = -72.933182332955
# --- Row 22 ---
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 3
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E22").Value = -100
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("G22").Value = 3
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = -50
# --- Row 23 ---
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -75
$ws.Range("F23").Value = 10
$ws.Range("H23").Value = -41.176470588235
$ws.Range("I23").Value = 89
$ws.Range("J23").Value = 99
$ws.Range("K23").Value = -10.101010101010
$ws.Range("L23").Value = 15.584415584415
$ws.Range("M23").Value = 53.448275862069
# --- Row 24 ---
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 62.5
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 33
$ws.Range("H24").Value = 15.151515151515
$ws.Range("I24").Value = 256
$ws.Range("J24").Value = 321
$ws.Range("K24").Value = -20.249221183800
$ws.Range("L24").Value = -20.249221183800
$ws.Range("M24").Value = -34.020618556701
# --- Row 25 ---
$ws.Range("I25").Value = 37
$ws.Range("J25").Value = 110
$ws.Range("K25").Value = -66.363636363636
$ws.Range("L25").Value = -65.420560747663
# --- Row 26 ---
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -83.333333333333
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = -73.684210526315
$ws.Range("I26").Value = 151
$ws.Range("J26").Value = 159
$ws.Range("K26").Value = -5.031446540880
$ws.Range("L26").Value = 6.338028169014
$ws.Range("M26").Value = -23.737373737373
# --- Row 27 ---
$ws.Range("D14").Copy($ws.Range("C27"))
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 42.857142857142
$ws.Range("L27").Value = -37.5
# --- Row 28 ---
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 24
$ws.Range("K28").Value = -11.111111111111
$ws.Range("L28").Value = -25
# --- Row 33 ---
$ws.Range("L33").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L33").Value = -100
